$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unused rows 6 and 7 (ECs/Resolving-Mac sending clusters removed)
$ws.Rows("6:7").Delete()

# Update data rows 2-5 with new TPM-derived values
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Rspo3"
$ws.Range("C2").Value = "Lgr6"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.583115666666667
$ws.Range("H2").Value = 22.749347
$ws.Range("I2").Value = 0.9996196502830235
$ws.Range("J2").Value = 0.9996196502830236
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05147733333333334
$ws.Range("N2").Value = 0.154432
$ws.Range("O2").Value = 0.887188413789934
$ws.Range("P2").Value = 0.8871884137899338
$ws.Range("Q2").Value = 0.3903585728782223
$ws.Range("R2").Value = 3.513227155904
$ws.Range("S2").Value = 0.8868509719278441
$ws.Range("T2").Value = 0.8868509719278441

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Rspo3"
$ws.Range("C3").Value = "Lgr6"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.583115666666667
$ws.Range("H3").Value = 22.749347
$ws.Range("I3").Value = 0.9996196502830235
$ws.Range("J3").Value = 0.9996196502830236
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.006545666666666668
$ws.Range("N3").Value = 0.019637
$ws.Range("O3").Value = 0.1128115862100661
$ws.Range("P3").Value = 0.1128115862100661
$ws.Range("Q3").Value = 0.04963654744877779
$ws.Range("R3").Value = 0.446728927039
$ws.Range("S3").Value = 0.1127686783551795
$ws.Range("T3").Value = 0.1127686783551795

$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Rspo3"
$ws.Range("C4").Value = "Lgr6"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.002885333333333333
$ws.Range("H4").Value = 0.008656
$ws.Range("I4").Value = 0.0003803497169764852
$ws.Range("J4").Value = 0.0003803497169764852
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05147733333333334
$ws.Range("N4").Value = 0.154432
$ws.Range("O4").Value = 0.887188413789934
$ws.Range("P4").Value = 0.8871884137899338
$ws.Range("Q4").Value = 0.0001485292657777778
$ws.Range("R4").Value = 0.001336763392
$ws.Range("S4").Value = 0.0003374418620898182
$ws.Range("T4").Value = 0.0003374418620898182

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Rspo3"
$ws.Range("C5").Value = "Lgr6"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.002885333333333333
$ws.Range("H5").Value = 0.008656
$ws.Range("I5").Value = 0.0003803497169764852
$ws.Range("J5").Value = 0.0003803497169764852
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.006545666666666668
$ws.Range("N5").Value = 0.019637
$ws.Range("O5").Value = 0.1128115862100661
$ws.Range("P5").Value = 0.1128115862100661
$ws.Range("Q5").Value = 0.00001888643022222223
$ws.Range("R5").Value = 0.000169977872
$ws.Range("S5").Value = 0.00004290785488666701
$ws.Range("T5").Value = 0.00004290785488666701
